$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("devices")

# Insert a new row above row 2, shifting existing data down.
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Clear()

$ws.Range("E2").Value = "Galaxy S5"
$ws.Range("I2").Value = "Brian"
$ws.Range("J2").Value = "mobileOS"

$ws.Range("I2").Select()
